# Hydrogen Production Efficiency by Pathway - add "guaranteed clean" pathway rows
# (electrolysis with guaranteed clean electricity / natural gas reforming with CCS)
# and drop the stray "+46" term from the natural-gas-reforming baseline formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPEbP")

# 1. Natural gas reforming "Today" value: 118/(162+2+46) -> 118/(162+2)
$ws.Range("B3").Formula = "=118/(162+2)"

# 2. New row 7: "electrolysis with guaranteed clean electricity" = copy of row 2
$ws.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$ws.Range("B7").Formula = "=B2"
$ws.Range("C7:AI7").Formula = "=C2"

# 3. New row 8: "natural gas reforming with CCS" = copy of row 3
$ws.Range("A8").Value = "natural gas reforming with CCS"
$ws.Range("B8").Formula = "=B3"
$ws.Range("C8:AI8").Formula = "=C3"

# 4. Re-apply the numeric format on the existing efficiency cells so the
#    redundant/duplicate cell style (a leftover, border-less clone of the
#    "0.000" style) collapses into the canonical one, same as Excel does
#    when it resaves this workbook.
$ws.Range("Q2:AI6").NumberFormat = "0.000"
